$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 344.22223
$ws.Range("I18").Value = 212.375
$ws.Range("J18").Value = 1399
$ws.Range("K18").Value = 212.375
$ws.Range("L18").Value = 1399
$ws.Range("M18").Value = 71.625
$ws.Range("N18").Value = -1967
$ws.Range("H137").Value = 10516.091
$ws.Range("I137").Value = 1245.5
$ws.Range("J137").Value = 12576.223
$ws.Range("K137").Value = 3736.5
$ws.Range("L137").Value = 37728.669
$ws.Range("M137").Value = -1186.5
$ws.Range("N137").Value = -42828.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 17131.125
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 17131.125
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 17131.125
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -18107.125
$ws.Range("H55").Value = 9000
$ws.Range("J55").Value = 9000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9630
$ws.Range("H63").Value = 2800.4546
$ws.Range("I63").Value = 2820.5
$ws.Range("J63").Value = 2600
$ws.Range("K63").Value = 2820.5
$ws.Range("L63").Value = 2600
$ws.Range("M63").Value = -2134.5
$ws.Range("N63").Value = -3972
$ws.Range("H66").Value = 2800.4546
$ws.Range("I66").Value = 2820.5
$ws.Range("J66").Value = 2600
$ws.Range("K66").Value = 14102.5
$ws.Range("L66").Value = 13000
$ws.Range("M66").Value = -10670.5
$ws.Range("N66").Value = -19864
$ws.Range("H80").Value = 18777.416
$ws.Range("J80").Value = 18777.416
$ws.Range("L80").Value = 18777.416
$ws.Range("N80").Value = -20773.416
$ws.Range("H83").Value = 18777.416
$ws.Range("J83").Value = 18777.416
$ws.Range("L83").Value = 56332.24800000001
$ws.Range("N83").Value = -66316.24800000001
$ws.Range("H88").Value = 14643.6875
$ws.Range("I88").Value = 1777.6666
$ws.Range("J88").Value = 31185.715
$ws.Range("K88").Value = 1777.6666
$ws.Range("L88").Value = 31185.715
$ws.Range("M88").Value = -1371.6666
$ws.Range("N88").Value = -31997.715
$ws.Range("H91").Value = 14643.6875
$ws.Range("I91").Value = 1777.6666
$ws.Range("J91").Value = 31185.715
$ws.Range("K91").Value = 1777.6666
$ws.Range("L91").Value = 31185.715
$ws.Range("M91").Value = -373.6666
$ws.Range("N91").Value = -33993.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 16160
$ws.Range("I35").Value = 5000
$ws.Range("J35").Value = 16957.143
$ws.Range("K35").Value = 5000
$ws.Range("L35").Value = 16957.143
$ws.Range("M35").Value = -4690
$ws.Range("N35").Value = -17577.143
$ws.Range("H82").Value = 11142.667
$ws.Range("I82").Value = 7371.4
$ws.Range("J82").Value = 29999
$ws.Range("K82").Value = 7371.4
$ws.Range("L82").Value = 29999
$ws.Range("M82").Value = -6988.4
$ws.Range("N82").Value = -30765
$ws.Range("H85").Value = 11142.667
$ws.Range("I85").Value = 7371.4
$ws.Range("J85").Value = 29999
$ws.Range("K85").Value = 7371.4
$ws.Range("L85").Value = 29999
$ws.Range("M85").Value = -6045.4
$ws.Range("N85").Value = -32651
$ws.Range("H86").Value = 320021.1
$ws.Range("I86").Value = 1930.3334
$ws.Range("J86").Value = 701730
$ws.Range("K86").Value = 1930.3334
$ws.Range("L86").Value = 701730
$ws.Range("M86").Value = -807.3334
$ws.Range("N86").Value = -703976
$ws.Range("H89").Value = 320021.1
$ws.Range("I89").Value = 1930.3334
$ws.Range("J89").Value = 701730
$ws.Range("K89").Value = 9651.666999999999
$ws.Range("L89").Value = 3508650
$ws.Range("M89").Value = -4035.666999999999
$ws.Range("N89").Value = -3519882
$ws.Range("H105").Value = 1475
$ws.Range("J105").Value = 1300
$ws.Range("L105").Value = 1300
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 61.266666
$ws.Range("J7").Value = 49.9
$ws.Range("L7").Value = 49.9
$ws.Range("N7").Value = -275.9
$ws.Range("H41").Value = 8538.666999999999
$ws.Range("J41").Value = 11200
$ws.Range("L41").Value = 11200
$ws.Range("N41").Value = -12056
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H59").Value = 12258.462
$ws.Range("J59").Value = 12258.462
$ws.Range("L59").Value = 12258.462
$ws.Range("N59").Value = -14548.462
$ws.Range("H68").Value = 26971.25
$ws.Range("I68").Value = 9000
$ws.Range("K68").Value = 9000
$ws.Range("M68").Value = -8251
$ws.Range("H71").Value = 26971.25
$ws.Range("I71").Value = 9000
$ws.Range("K71").Value = 27000
$ws.Range("M71").Value = -23256
$ws.Range("H74").Value = 12577.818
$ws.Range("J74").Value = 12577.818
$ws.Range("L74").Value = 12577.818
$ws.Range("N74").Value = -14325.818
$ws.Range("H77").Value = 12577.818
$ws.Range("J77").Value = 12577.818
$ws.Range("L77").Value = 37733.454
$ws.Range("N77").Value = -46469.454
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 38618612
$ws.Range("I131").Value = 390
$ws.Range("J131").Value = 54598564
$ws.Range("K131").Value = 1170
$ws.Range("L131").Value = 163795692
$ws.Range("M131").Value = 3870
$ws.Range("N131").Value = -163805772
$ws.Range("H140").Value = 4837386.5
$ws.Range("I140").Value = 1748.25
$ws.Range("K140").Value = 5244.75
$ws.Range("M140").Value = -64.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 17000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -17224
$ws.Range("H70").Value = 4337.5415
$ws.Range("I70").Value = 4292.857
$ws.Range("K70").Value = 4292.857
$ws.Range("M70").Value = -4022.857
$ws.Range("H73").Value = 4337.5415
$ws.Range("I73").Value = 4292.857
$ws.Range("K73").Value = 4292.857
$ws.Range("M73").Value = -3356.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 51684.668
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 62001.6
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 62001.6
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = -62225.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2673666.8
$ws.Range("J2").Value = 10500
$ws.Range("L2").Value = 10500
$ws.Range("N2").Value = -10724
$ws.Range("H136").Value = 2395811.5
$ws.Range("I136").Value = 2463949.2
$ws.Range("J136").Value = 2000611.8
$ws.Range("K136").Value = 7391847.600000001
$ws.Range("L136").Value = 6001835.4
$ws.Range("M136").Value = -7389297.600000001
$ws.Range("N136").Value = -6006935.4
